# Auto-generated from the provided OOXML diff.
# Each target cell keeps its original "text" semantics (as in before.xlsx,
# every data cell is an inline/shared string, even values that look numeric,
# e.g. "1.00" or "505.78"). Using a leading apostrophe forces Excel/IronCalc
# to store the new value as text instead of silently converting it to a
# number, and resetting .Style to "Normal" afterwards clears the transient
# quote-prefix flag so the cell keeps the same (default) style as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue 'D2' '59.999.10'
Set-TextValue 'E2' '  +0.53%  '
Set-TextValue 'D3' '2.566.31'
Set-TextValue 'E3' '  -0.51%  '
Set-TextValue 'E4' '  +0.15%  '
Set-TextValue 'D5' '505.21'
Set-TextValue 'E5' '  -0.03%  '
Set-TextValue 'D6' '152.14'
Set-TextValue 'E6' '  -2.98%  '
Set-TextValue 'D7' '0.999'
Set-TextValue 'E7' '  +0.29%  '
Set-TextValue 'E8' '  -5.28%  '
Set-TextValue 'D9' '2.569.84'
Set-TextValue 'E9' '  -0.32%  '
Set-TextValue 'D10' '6.53'
Set-TextValue 'E10' '  +7.22%  '
Set-TextValue 'D11' '0.103'
Set-TextValue 'E11' '  +0.36%  '
Set-TextValue 'D12' '0.347'
Set-TextValue 'E12' '  +1.75%  '
Set-TextValue 'E13' '  +1.13%  '
Set-TextValue 'D14' '3.020.28'
Set-TextValue 'E14' '  +0.76%  '
Set-TextValue 'D15' '60.101.66'
Set-TextValue 'E15' '  +1.20%  '
Set-TextValue 'D16' '21.45'
Set-TextValue 'E16' '  -1.72%  '
Set-TextValue 'D17' '0.0000139'
Set-TextValue 'E17' '  +1.58%  '
Set-TextValue 'D18' '2.571.84'
Set-TextValue 'E18' '  -0.04%  '
Set-TextValue 'D19' '4.76'
Set-TextValue 'E19' '  +0.12%  '
Set-TextValue 'D20' '344.52'
Set-TextValue 'E20' '  +1.40%  '
Set-TextValue 'D21' '10.38'
Set-TextValue 'E21' '  +0.25%  '
Set-TextValue 'D22' '6.08'
Set-TextValue 'E22' '  +0.65%  '
Set-TextValue 'D23' '0.998'
Set-TextValue 'E23' '  -0.34%  '
Set-TextValue 'D24' '59.70'
Set-TextValue 'E24' '  -0.47%  '
Set-TextValue 'D25' '0.418'
Set-TextValue 'E25' '  +0.09%  '
Set-TextValue 'E26' '  +0.11%  '
Set-TextValue 'D27' '1.00'
Set-TextValue 'E27' '  +0.43%  '
Set-TextValue 'D28' '0.0₃0836'
Set-TextValue 'E28' '  +0.29%  '
Set-TextValue 'D29' '7.34'
Set-TextValue 'E29' '  -0.01%  '
Set-TextValue 'E30' '  +0.34%  '
Set-TextValue 'D31' '19.24'
Set-TextValue 'E31' '  -0.98%  '
Set-TextValue 'D32' '153.08'
Set-TextValue 'E32' '  -2.58%  '
Set-TextValue 'E33' '  -1.09%  '
Set-TextValue 'D34' '5.69'
Set-TextValue 'E34' '  +3.18%  '
Set-TextValue 'E35' '  +1.44%  '
Set-TextValue 'D36' '1.18'
Set-TextValue 'E36' '  -1.33%  '
Set-TextValue 'D37' '0.848'
Set-TextValue 'E37' '  +7.92%  '
Set-TextValue 'D38' '0.846'
Set-TextValue 'E38' '  -1.45%  '
Set-TextValue 'E39' '  +1.70%  '
Set-TextValue 'B40' 'OKB'
Set-TextValue 'C40' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D40' '35.97'
Set-TextValue 'E40' '  +2.15%  '
Set-TextValue 'B41' 'Filecoin'
Set-TextValue 'C41' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D41' '3.73'
Set-TextValue 'E41' '  -0.13%  '
Set-TextValue 'D42' '292.32'
Set-TextValue 'E42' '  -2.99%  '
Set-TextValue 'D43' '0.617'
Set-TextValue 'E43' '  -2.12%  '
Set-TextValue 'D44' '0.0991'
Set-TextValue 'E44' '  -2.43%  '
Set-TextValue 'D45' '0.998'
Set-TextValue 'E45' '  -0.19%  '
Set-TextValue 'D46' '0.0553'
Set-TextValue 'E46' '  -3.89%  '
Set-TextValue 'D47' '19.66'
Set-TextValue 'E47' '  +2.29%  '
Set-TextValue 'D48' '4.83'
Set-TextValue 'E48' '  -1.83%  '
Set-TextValue 'E49' '  -2.05%  '
Set-TextValue 'E50' '  +0.35%  '
Set-TextValue 'D51' '1.991.09'
Set-TextValue 'E51' '  +0.31%  '
